# Actualización automática 2025-06-01 08:00:06
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": zero-out this period's PORCELANATO (L) figures
# plus the two stray D18/E6/F6 values and refresh the "x de 17" counters.
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("L2").Value = 0
$wsGrupo.Range("L3").Value = 0
$wsGrupo.Range("L5").Value = 0
$wsGrupo.Range("E6").Value = 0
$wsGrupo.Range("F6").Value = 0
$wsGrupo.Range("L6").Value = 0
$wsGrupo.Range("L8").Value = 0
$wsGrupo.Range("L10").Value = 0
$wsGrupo.Range("L14").Value = 0
$wsGrupo.Range("L15").Value = 0
$wsGrupo.Range("D18").Value = 0
$wsGrupo.Range("L18").Value = 0

$wsGrupo.Range("D19").Value = "0 de 17"
$wsGrupo.Range("E19").Value = "0 de 17"
$wsGrupo.Range("F19").Value = "0 de 17"
$wsGrupo.Range("L19").Value = "0 de 17"

# --- Sheet "VENTA MENSUAL": roll the monthly columns forward one month
# (febrero drops off, junio is appended) by shifting each client's figures
# left and relabelling the header row.
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("C1").Value = "marzo"
$wsMensual.Range("D1").Value = "abril"
$wsMensual.Range("E1").Value = "mayo"
$wsMensual.Range("F1").Value = "junio"

for ($r = 2; $r -le 19; $r++) {
    $oldD = $wsMensual.Cells.Item($r, 4).Value2
    $oldE = $wsMensual.Cells.Item($r, 5).Value2
    $oldF = $wsMensual.Cells.Item($r, 6).Value2

    $wsMensual.Cells.Item($r, 3).Value = $oldD
    $wsMensual.Cells.Item($r, 4).Value = $oldE
    $wsMensual.Cells.Item($r, 5).Value = $oldF
    $wsMensual.Cells.Item($r, 6).Value = 0
}

# NOTE: the host's ColumnWidth setter round-trips through pixel units and
# always adds back 5/6 of a character, so back that off here to land on the
# exact widths (14 / 11) recorded in the sheet's <cols>.
$wsMensual.Columns.Item(4).ColumnWidth = 14 - (5 / 6)
$wsMensual.Columns.Item(6).ColumnWidth = 11 - (5 / 6)
